$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-10-30 04:27:31"

foreach ($sheetName in @("Главные", "Линейные")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA$row").Value = $newTimestamp
    }
}
